$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before the "Late" column (old column N),
# shifting Late/heading(Date)/Outstanding one column to the right.
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = 10.7109375

# Make "Repayment schedule" the active/selected sheet with the new selection.
$ws.Activate()
$ws.Range("S7").Select()
